$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match E7's cell style (s="1") to the one already used on E4, before
# filling in values, by copying formats from E4
$ws.Range("E4").Copy()
$ws.Range("E7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the remaining cells of row 7 (C7, D7, E7, F7)
$ws.Range("C7").Value = 0.54861111111111105
$ws.Range("D7").Value = "HARD"
$ws.Range("E7").Value = "Aula de HARD (JWT) +  implemetação produtos do projeto green collections"
$ws.Range("F7").Value = 0.010416666666666666

# Row 7 grows taller (wrapped subject text), matching the other data rows
$ws.Range("A7:F7").RowHeight = 31.5

# Update the view: scroll so column B is the left-most visible column,
# and move the active selection to E8
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E8").Select()
